# Add a "Letter" column (E) that records the document type ("PDF") for
# every reference response row, so references can be matched to students
# by how they submitted their letter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row (header + one row per reference response).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# New header in E1, matching the other header cells.
$ws.Range("E1").Value = "Letter"

# Every data row (2..lastRow) records "PDF" as the letter format.
$dataRange = $ws.Range("E2:E" + $lastRow)
$dataRange.Value = "PDF"

# Copy A1's cell format (font/fill/alignment) onto the new column so it
# matches the style used for the rest of the header/data styling.
$ws.Range("A1").Copy()
$ws.Range("E1:E" + $lastRow).PasteSpecial(-4122)
